# mozzarella db params — fix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide the "working" columns E:U (5-21); keep the existing custom
# widths on M (13) and Q (17), everything else keeps its (default) width.
$ws.Range("E1:U1").EntireColumn.Hidden = $true

# --- V column: 50 -> 35 for this set of rows.
$rowsToFix = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,37,38,41,42,43,44,49,50,51,52,56,65,66,70,73,74,77,78)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 22).Value2 = 35
}

# --- Row height tweaks.
$rowsToResize = @(20,56,59,60,61)
foreach ($r in $rowsToResize) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# --- View: scroll position + active cell/selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("W24").Select()
